# Update the "Hypertension" (HT) method resource workbook:
#  - refresh age-banded prevalence figures on prevalence2018 (column C, rows 20-122)
#    the incidence2018_plus sheet pulls these via =prevalence2018!C.../100 so it
#    recalculates automatically.
#  - leave the UI in the state the author left it in: "data" sheet selection
#    parked at E26, "prevalence2018" the active tab, scrolled/selected near H51.

$wb = $excel.ActiveWorkbook

# --- 1) prevalence2018: new age-banded prevalence values ------------------
$ws = $wb.Worksheets.Item("prevalence2018")

$ws.Range("C20").Value = 0
$ws.Range("C21:C36").Value = 0.35
$ws.Range("C37:C46").Value = 0.43
$ws.Range("C47:C56").Value = 0.56999999999999995
$ws.Range("C57:C122").Value = 0.9

# --- 2) Window/selection state -------------------------------------------
# Park the selection on "data" (previously the active tab) at E26 before
# switching away, so its stored selection matches the author's last action.
$wsData = $wb.Worksheets.Item("data")
$wsData.Activate()
$wsData.Range("E26").Select()

# Make "prevalence2018" the active sheet again, scrolled down near row 47
# with H51 selected (matches the final on-screen state of the edit).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 47
$ws.Range("H51").Select()

Write-Output "Updated HT prevalence figures"
